$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historias de Usuario")

# --- Re-shuffle the B:E content (and formatting) of rows 3-7 -------------
# The HU id (column A) stays put on each row, but the Descripcion /
# Criterios Aceptacion / Tamano columns move between rows:
#   new row3 <- old row4
#   new row4 <- old row5
#   new row5 <- old row7
#   new row6 <- old row3
#   new row7 <- old row6
# Stage each source block in a scratch area first so the cyclic
# permutation doesn't clobber a source before it has been read.
$ws.Range("B3:E3").Copy($ws.Range("B100:E100"))
$ws.Range("B4:E4").Copy($ws.Range("B101:E101"))
$ws.Range("B5:E5").Copy($ws.Range("B102:E102"))
$ws.Range("B6:E6").Copy($ws.Range("B103:E103"))
$ws.Range("B7:E7").Copy($ws.Range("B104:E104"))

$ws.Range("B101:E101").Copy($ws.Range("B3:E3"))
$ws.Range("B102:E102").Copy($ws.Range("B4:E4"))
$ws.Range("B104:E104").Copy($ws.Range("B5:E5"))
$ws.Range("B100:E100").Copy($ws.Range("B6:E6"))
$ws.Range("B103:E103").Copy($ws.Range("B7:E7"))

$ws.Range("B100:E104").ClearContents()

# --- Row 4 grew taller to fit its (new) content ---------------------------
$ws.Rows.Item(4).RowHeight = 250.5

# --- View state: frozen pane now shows the header row on top, selection
#     moved back up to B1 -------------------------------------------------
$ws.Activate()
$ws.Range("B1").Select()
